# --------------------------------------------------------------------------
# Edit 1: drop the lone "test" run (paragraph w14:paraId="0B622489").
# "test" occurs exactly once in the whole document body, so a single
# Find/Replace-with-nothing removes that run cleanly (the now-empty <w:r>
# is dropped, matching the diff exactly).
# --------------------------------------------------------------------------
$d = $word.ActiveDocument
$null = $d.Content.Find.Execute("test", $true, $true, $false, $false, $false, `
                                 $true, 1, $false, "", 2)

# --------------------------------------------------------------------------
# Edit 2: fill in the empty paragraph (w14:paraId="4620C906") that sits right
# after the "Rapportens kapitel" heading and right before "Slutsatser och
# avslutande diskussion" with the new "Spelregler" / "Roll och
# ansvarsfordelning content block -- several new paragraphs, each with
# its own run/paragraph formatting (theme fonts, bold, sizes, manual breaks).
# We build the exact OOXML for the replacement and hand it to the paragraph's
# Range via InsertXML (WordOpenXML package) so formatting round-trips exactly.
# --------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "`r") {
        if ($i -gt 1 -and $i -lt $d.Paragraphs.Count) {
            $prevText = $d.Paragraphs.Item($i - 1).Range.Text
            $nextText = $d.Paragraphs.Item($i + 1).Range.Text
            if ($prevText -like "*Rapportens kapitel*" -and $nextText -like "*Slutsatser*") {
                $targetIndex = $i
                break
            }
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the target empty paragraph after 'Rapportens kapitel'"
}

$target = $d.Paragraphs.Item($targetIndex).Range

$paraSpelreglerKommunikation = '<w:p w14:paraId="4620C906" w14:textId="77777777" w:rsidR="00D37610" w:rsidRDefault="00D37610" w:rsidP="00832785"><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Spelregler</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Vad gäller spelregler är det bestämt att vi sitter från 9-16 varje dag i samtal. 16:30 görs dagligen en sammanfattning på dagen och vad som ska göras inför nästkommande dag.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/><w:t xml:space="preserve">Sedan tidigare har gruppen skapat ett set med spelregler </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>som även är tänkta att efterföljas i detta arbete. Följande punkter togs då upp:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Kommunikation</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Om en person i gruppen får svårigheter med sina uppgifter är det bra att kommunicera det till gruppen i tid.</w:t></w:r></w:p>'

$paraHallaTider = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Hålla tider/anmäla frånvaro</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ha tålamod och förståelse för formatet vi jobbar i</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ha förståelse för att kommunikation inte är densamma vid ett fysiskt möte</w:t></w:r></w:p>'

$paraStruktur = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Struktur</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Att arbeta strukturerat och veta hur vi gör och vad vi gör, annars får vi stämma av med varandra</w:t></w:r></w:p>'

$paraSamarbeta = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Samarbeta</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Dela information med varandra</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, förutsatt att detta sker inom rimliga gränser.</w:t></w:r></w:p>'

$paraOmViMarker = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Om vi märker att dessa inte går att upprätthålla så är det viktigt att i tidigt skede kommunicera detta och agera enligt Newtons riktlinjer för projektarbete.</w:t></w:r></w:p>'

$paraBlankSpacer = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/></w:pPr></w:p>'

$paraRollOchAnsvar = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7938"/></w:tabs><w:ind w:right="992"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Roll och ansvarsfördelning</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:t>Vad gäller roller och ansvarsfördelning</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Scrum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>Nedbrytning av projektet</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:t>För att få en överskådlig bild av omfattningen av projektet är det viktigt att dela upp projektet i mindre delar</w:t></w:r><w:r><w:t xml:space="preserve">. Vi har valt att dela upp projektet dels i de 2 olika releaserna som kommer att göras. Vi delar även in arbetet i dels en del som rör rapportskrivning </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>och planering</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>WBS</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Tidsplanering</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:t xml:space="preserve">Projektet pågår under 3 veckors tid. Efter den WBS som har gjorts </w:t></w:r></w:p>'

$newBodyContent = $paraSpelreglerKommunikation + $paraHallaTider + $paraStruktur + $paraSamarbeta + $paraOmViMarker + $paraBlankSpacer + $paraRollOchAnsvar

$xmlHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$xml = $xmlHeader + $newBodyContent + $xmlFooter

$target.InsertXML($xml)

Write-Output "OK targetIndex=$targetIndex"
